$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 26.731658
$ws.Range("H2").Value = 80.194974
$ws.Range("I2").Value = 0.02353393228912
$ws.Range("J2").Value = 0.02353393228912
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.043534
$ws.Range("N2").Value = 0.130602
$ws.Range("O2").Value = 0.00760428821720775
$ws.Range("P2").Value = 0.007604288217207752
$ws.Range("Q2").Value = 1.163735999372
$ws.Range("R2").Value = 10.473623994348
$ws.Range("S2").Value = 0.0001789588040107202
$ws.Range("T2").Value = 0.0001789588040107203
$ws.Range("G3").Value = 26.731658
$ws.Range("H3").Value = 80.194974
$ws.Range("I3").Value = 0.02353393228912
$ws.Range("J3").Value = 0.02353393228912
$ws.Range("O3").Value = 0.964981158713912
$ws.Range("P3").Value = 0.9649811587139122
$ws.Range("Q3").Value = 147.6776367536787
$ws.Range("R3").Value = 1329.098730783108
$ws.Range("S3").Value = 0.02270980124944976
$ws.Range("T3").Value = 0.02270980124944976
$ws.Range("G4").Value = 26.731658
$ws.Range("H4").Value = 80.194974
$ws.Range("I4").Value = 0.02353393228912
$ws.Range("J4").Value = 0.02353393228912
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03706533333333333
$ws.Range("N4").Value = 0.111196
$ws.Range("O4").Value = 0.006474375833453032
$ws.Range("P4").Value = 0.006474375833453035
$ws.Range("Q4").Value = 0.9908178143226666
$ws.Range("R4").Value = 8.917360328904001
$ws.Range("S4").Value = 0.0001523675224787985
$ws.Range("T4").Value = 0.0001523675224787986
$ws.Range("G5").Value = 26.731658
$ws.Range("H5").Value = 80.194974
$ws.Range("I5").Value = 0.02353393228912
$ws.Range("J5").Value = 0.02353393228912
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.119881
$ws.Range("N5").Value = 0.359643
$ws.Range("O5").Value = 0.02094017723542708
$ws.Range("P5").Value = 0.02094017723542708
$ws.Range("Q5").Value = 3.204617892698
$ws.Range("R5").Value = 28.84156103428201
$ws.Range("S5").Value = 0.0004928047131807128
$ws.Range("T5").Value = 0.000492804713180713
$ws.Range("I6").Value = 0.9376016087099961
$ws.Range("J6").Value = 0.9376016087099961
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.043534
$ws.Range("N6").Value = 0.130602
$ws.Range("O6").Value = 0.00760428821720775
$ws.Range("P6").Value = 0.007604288217207752
$ws.Range("Q6").Value = 46.36372416306133
$ws.Range("R6").Value = 417.273517467552
$ws.Range("S6").Value = 0.007129792865548454
$ws.Range("T6").Value = 0.007129792865548457
$ws.Range("I7").Value = 0.9376016087099961
$ws.Range("J7").Value = 0.9376016087099961
$ws.Range("O7").Value = 0.964981158713912
$ws.Range("P7").Value = 0.9649811587139122
$ws.Range("S7").Value = 0.904767886785
$ws.Range("T7").Value = 0.9047678867850002
$ws.Range("I8").Value = 0.9376016087099961
$ws.Range("J8").Value = 0.9376016087099961
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.03706533333333333
$ws.Range("N8").Value = 0.111196
$ws.Range("O8").Value = 0.006474375833453032
$ws.Range("P8").Value = 0.006474375833453035
$ws.Range("Q8").Value = 39.47459205858845
$ws.Range("R8").Value = 355.271328527296
$ws.Range("S8").Value = 0.006070385196838685
$ws.Range("T8").Value = 0.006070385196838687
$ws.Range("I9").Value = 0.9376016087099961
$ws.Range("J9").Value = 0.9376016087099961
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.119881
$ws.Range("N9").Value = 0.359643
$ws.Range("O9").Value = 0.02094017723542708
$ws.Range("P9").Value = 0.02094017723542708
$ws.Range("Q9").Value = 127.6733040012854
$ws.Range("R9").Value = 1149.059736011568
$ws.Range("S9").Value = 0.01963354386260887
$ws.Range("T9").Value = 0.01963354386260887
$ws.Range("G10").Value = 0.1721486666666666
$ws.Range("H10").Value = 0.516446
$ws.Range("I10").Value = 0.0001515556971810586
$ws.Range("J10").Value = 0.0001515556971810586
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.043534
$ws.Range("N10").Value = 0.130602
$ws.Range("O10").Value = 0.00760428821720775
$ws.Range("P10").Value = 0.007604288217207752
$ws.Range("Q10").Value = 0.007494320054666665
$ws.Range("R10").Value = 0.067448880492
$ws.Range("S10").Value = 0.000001152473202324629
$ws.Range("T10").Value = 0.00000115247320232463
$ws.Range("G11").Value = 0.1721486666666666
$ws.Range("H11").Value = 0.516446
$ws.Range("I11").Value = 0.0001515556971810586
$ws.Range("J11").Value = 0.0001515556971810586
$ws.Range("O11").Value = 0.964981158713912
$ws.Range("P11").Value = 0.9649811587139122
$ws.Range("Q11").Value = 0.9510262425035555
$ws.Range("R11").Value = 8.559236182531999
$ws.Range("S11").Value = 0.0001462483922754727
$ws.Range("T11").Value = 0.0001462483922754727
$ws.Range("G12").Value = 0.1721486666666666
$ws.Range("H12").Value = 0.516446
$ws.Range("I12").Value = 0.0001515556971810586
$ws.Range("J12").Value = 0.0001515556971810586
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.03706533333333333
$ws.Range("N12").Value = 0.111196
$ws.Range("O12").Value = 0.006474375833453032
$ws.Range("P12").Value = 0.006474375833453035
$ws.Range("Q12").Value = 0.006380747712888888
$ws.Range("R12").Value = 0.05742672941599999
$ws.Range("S12").Value = 0.0000009812285432511714
$ws.Range("T12").Value = 0.0000009812285432511719
$ws.Range("G13").Value = 0.1721486666666666
$ws.Range("H13").Value = 0.516446
$ws.Range("I13").Value = 0.0001515556971810586
$ws.Range("J13").Value = 0.0001515556971810586
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.119881
$ws.Range("N13").Value = 0.359643
$ws.Range("O13").Value = 0.02094017723542708
$ws.Range("P13").Value = 0.02094017723542708
$ws.Range("Q13").Value = 0.02063735430866667
$ws.Range("R13").Value = 0.185736188778
$ws.Range("S13").Value = 0.000003173603160010082
$ws.Range("T13").Value = 0.000003173603160010084
$ws.Range("G14").Value = 43.33877
$ws.Range("H14").Value = 130.01631
$ws.Range("I14").Value = 0.03815444888131313
$ws.Range("J14").Value = 0.03815444888131313
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.043534
$ws.Range("N14").Value = 0.130602
$ws.Range("O14").Value = 0.00760428821720775
$ws.Range("P14").Value = 0.007604288217207752
$ws.Range("Q14").Value = 1.88671001318
$ws.Range("R14").Value = 16.98039011862
$ws.Range("S14").Value = 0.0002901374260622248
$ws.Range("T14").Value = 0.0002901374260622249
$ws.Range("G15").Value = 43.33877
$ws.Range("H15").Value = 130.01631
$ws.Range("I15").Value = 0.03815444888131313
$ws.Range("J15").Value = 0.03815444888131313
$ws.Range("O15").Value = 0.964981158713912
$ws.Range("P15").Value = 0.9649811587139122
$ws.Range("Q15").Value = 239.4227523564467
$ws.Range("R15").Value = 2154.80477120802
$ws.Range("S15").Value = 0.03681832429158027
$ws.Range("T15").Value = 0.03681832429158027
$ws.Range("G16").Value = 43.33877
$ws.Range("H16").Value = 130.01631
$ws.Range("I16").Value = 0.03815444888131313
$ws.Range("J16").Value = 0.03815444888131313
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.03706533333333333
$ws.Range("N16").Value = 0.111196
$ws.Range("O16").Value = 0.006474375833453032
$ws.Range("P16").Value = 0.006474375833453035
$ws.Range("Q16").Value = 1.606365956306667
$ws.Range("R16").Value = 14.45729360676
$ws.Range("S16").Value = 0.0002470262417758928
$ws.Range("T16").Value = 0.0002470262417758929
$ws.Range("G17").Value = 43.33877
$ws.Range("H17").Value = 130.01631
$ws.Range("I17").Value = 0.03815444888131313
$ws.Range("J17").Value = 0.03815444888131313
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.119881
$ws.Range("N17").Value = 0.359643
$ws.Range("O17").Value = 0.02094017723542708
$ws.Range("P17").Value = 0.02094017723542708
$ws.Range("Q17").Value = 5.195495086370001
$ws.Range("R17").Value = 46.75945577733001
$ws.Range("S17").Value = 0.0007989609218947393
$ws.Range("T17").Value = 0.0007989609218947395
$ws.Range("G18").Value = 0.4290093333333333
$ws.Range("H18").Value = 1.287028
$ws.Range("I18").Value = 0.0003776898762533613
$ws.Range("J18").Value = 0.0003776898762533613
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = 0.3333333333333333
$ws.Range("M18").Value = 0.043534
$ws.Range("N18").Value = 0.130602
$ws.Range("O18").Value = 0.00760428821720775
$ws.Range("P18").Value = 0.007604288217207752
$ws.Range("Q18").Value = 0.01867649231733333
$ws.Range("R18").Value = 0.168088430856
$ws.Range("S18").Value = 0.000002872062675752089
$ws.Range("T18").Value = 0.000002872062675752089
$ws.Range("G19").Value = 0.4290093333333333
$ws.Range("H19").Value = 1.287028
$ws.Range("I19").Value = 0.0003776898762533613
$ws.Range("J19").Value = 0.0003776898762533613
$ws.Range("O19").Value = 0.964981158713912
$ws.Range("P19").Value = 0.9649811587139122
$ws.Range("Q19").Value = 2.370039467508445
$ws.Range("R19").Value = 21.330355207576
$ws.Range("S19").Value = 0.0003644636144214826
$ws.Range("T19").Value = 0.0003644636144214827
$ws.Range("G20").Value = 0.4290093333333333
$ws.Range("H20").Value = 1.287028
$ws.Range("I20").Value = 0.0003776898762533613
$ws.Range("J20").Value = 0.0003776898762533613
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.3333333333333333
$ws.Range("M20").Value = 0.03706533333333333
$ws.Range("N20").Value = 0.111196
$ws.Range("O20").Value = 0.006474375833453032
$ws.Range("P20").Value = 0.006474375833453035
$ws.Range("Q20").Value = 0.01590137394311111
$ws.Range("R20").Value = 0.143112365488
$ws.Range("S20").Value = 0.000002445306207354629
$ws.Range("T20").Value = 0.00000244530620735463
$ws.Range("G21").Value = 0.4290093333333333
$ws.Range("H21").Value = 1.287028
$ws.Range("I21").Value = 0.0003776898762533613
$ws.Range("J21").Value = 0.0003776898762533613
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.119881
$ws.Range("N21").Value = 0.359643
$ws.Range("O21").Value = 0.02094017723542708
$ws.Range("P21").Value = 0.02094017723542708
$ws.Range("Q21").Value = 0.05143006788933334
$ws.Range("R21").Value = 0.462870611004
$ws.Range("S21").Value = 0.000007908892948771907
$ws.Range("T21").Value = 0.000007908892948771909
$ws.Range("G22").Value = 0.2053263333333333
$ws.Range("H22").Value = 0.6159789999999999
$ws.Range("I22").Value = 0.0001807645461362684
$ws.Range("J22").Value = 0.0001807645461362684
$ws.Range("K22").Value = 1
$ws.Range("L22").Value = 0.3333333333333333
$ws.Range("M22").Value = 0.043534
$ws.Range("N22").Value = 0.130602
$ws.Range("O22").Value = 0.00760428821720775
$ws.Range("P22").Value = 0.007604288217207752
$ws.Range("Q22").Value = 0.008938676595333331
$ws.Range("R22").Value = 0.08044808935799999
$ws.Range("S22").Value = 0.000001374585708272932
$ws.Range("T22").Value = 0.000001374585708272933
$ws.Range("G23").Value = 0.2053263333333333
$ws.Range("H23").Value = 0.6159789999999999
$ws.Range("I23").Value = 0.0001807645461362684
$ws.Range("J23").Value = 0.0001807645461362684
$ws.Range("O23").Value = 0.964981158713912
$ws.Range("P23").Value = 0.9649811587139122
$ws.Range("Q23").Value = 1.134314514646444
$ws.Range("R23").Value = 10.208830631818
$ws.Range("S23").Value = 0.0001744343811849707
$ws.Range("T23").Value = 0.0001744343811849707
$ws.Range("G24").Value = 0.2053263333333333
$ws.Range("H24").Value = 0.6159789999999999
$ws.Range("I24").Value = 0.0001807645461362684
$ws.Range("J24").Value = 0.0001807645461362684
$ws.Range("K24").Value = 1
$ws.Range("L24").Value = 0.3333333333333333
$ws.Range("M24").Value = 0.03706533333333333
$ws.Range("N24").Value = 0.111196
$ws.Range("O24").Value = 0.006474375833453032
$ws.Range("P24").Value = 0.006474375833453035
$ws.Range("Q24").Value = 0.00761048898711111
$ws.Range("R24").Value = 0.068494400884
$ws.Range("S24").Value = 0.000001170337609049762
$ws.Range("T24").Value = 0.000001170337609049762
$ws.Range("G25").Value = 0.2053263333333333
$ws.Range("H25").Value = 0.6159789999999999
$ws.Range("I25").Value = 0.0001807645461362684
$ws.Range("J25").Value = 0.0001807645461362684
$ws.Range("K25").Value = 3
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = 0.119881
$ws.Range("N25").Value = 0.359643
$ws.Range("O25").Value = 0.02094017723542708
$ws.Range("P25").Value = 0.02094017723542708
$ws.Range("Q25").Value = 0.02461472616633333
$ws.Range("R25").Value = 0.221532535497
$ws.Range("S25").Value = 0.000003785241633974995
$ws.Range("T25").Value = 0.000003785241633974997
